# Update "想去人数" (number of people interested) figures (column F)
# across the workbook's sheets, per the latest scrape refresh.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2905
$ws1.Range("F3").Value = 21346
$ws1.Range("F5").Value = 3200
$ws1.Range("F6").Value = 825
$ws1.Range("F8").Value = 535
$ws1.Range("F9").Value = 788
$ws1.Range("F10").Value = 294
$ws1.Range("F13").Value = 131
$ws1.Range("F16").Value = 301
$ws1.Range("F17").Value = 26
$ws1.Range("F18").Value = 438
$ws1.Range("F20").Value = 34
$ws1.Range("F22").Value = 50

# --- 演出 (Performances) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 164

# --- 本地生活 (Local life) sheet ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6150
$ws3.Range("F3").Value = 716
$ws3.Range("F4").Value = 711
$ws3.Range("F5").Value = 1682

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6150
$ws4.Range("F3").Value = 716
$ws4.Range("F4").Value = 711
$ws4.Range("F5").Value = 1682
$ws4.Range("F6").Value = 2905
$ws4.Range("F7").Value = 21346
$ws4.Range("F12").Value = 3200
$ws4.Range("F13").Value = 825
$ws4.Range("F17").Value = 535
$ws4.Range("F18").Value = 788
$ws4.Range("F19").Value = 294
$ws4.Range("F25").Value = 131
$ws4.Range("F28").Value = 540
$ws4.Range("F32").Value = 301
$ws4.Range("F33").Value = 164
$ws4.Range("F34").Value = 164
$ws4.Range("F35").Value = 26
$ws4.Range("F36").Value = 438
$ws4.Range("F39").Value = 34
$ws4.Range("F43").Value = 50
